$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01309533333333333
$ws.Range("H2").Value = 0.039286
$ws.Range("I2").Value = 0.000097606391451160887866127236
$ws.Range("J2").Value = 0.00009760639145116087431360008
$ws.Range("M2").Value = 3.442633333333333
$ws.Range("N2").Value = 10.3279
$ws.Range("O2").Value = 0.2185778494939793
$ws.Range("P2").Value = 0.2185778494939793
$ws.Range("Q2").Value = 0.04508243104444445
$ws.Range("R2").Value = 0.4057418794
$ws.Range("S2").Value = 0.000021334595140262269769796527
$ws.Range("T2").Value = 0.000021334595140262269769796527

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01309533333333333
$ws.Range("H3").Value = 0.039286
$ws.Range("I3").Value = 0.000097606391451160887866127236
$ws.Range("J3").Value = 0.00009760639145116087431360008
$ws.Range("O3").Value = 0.1871393710983698
$ws.Range("P3").Value = 0.1871393710983698
$ws.Range("Q3").Value = 0.03859813706088889
$ws.Range("R3").Value = 0.347383233548
$ws.Range("S3").Value = 0.000018265998711351549366535552
$ws.Range("T3").Value = 0.000018265998711351549366535552

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01309533333333333
$ws.Range("H4").Value = 0.039286
$ws.Range("I4").Value = 0.000097606391451160887866127236
$ws.Range("J4").Value = 0.00009760639145116087431360008
$ws.Range("M4").Value = 9.360041333333335
$ws.Range("O4").Value = 0.594282779407651
$ws.Range("P4").Value = 0.5942827794076508
$ws.Range("Q4").Value = 0.1225728612737778
$ws.Range("R4").Value = 1.103155751464
$ws.Range("S4").Value = 0.000058005797599547082282322313
$ws.Range("T4").Value = 0.000058005797599547061953531579

# Row 5
$ws.Range("I5").Value = 0.9877039077543427
$ws.Range("J5").Value = 0.9877039077543425
$ws.Range("M5").Value = 3.442633333333333
$ws.Range("N5").Value = 10.3279
$ws.Range("O5").Value = 0.2185778494939793
$ws.Range("P5").Value = 0.2185778494939793
$ws.Range("Q5").Value = 456.2005894454555
$ws.Range("R5").Value = 4105.8053050091
$ws.Range("S5").Value = 0.215890196093744
$ws.Range("T5").Value = 0.2158901960937439

# Row 6
$ws.Range("I6").Value = 0.9877039077543427
$ws.Range("J6").Value = 0.9877039077543425
$ws.Range("O6").Value = 0.1871393710983698
$ws.Range("P6").Value = 0.1871393710983698
$ws.Range("S6").Value = 0.18483828812855
$ws.Range("T6").Value = 0.1848382881285499

# Row 7
$ws.Range("I7").Value = 0.9877039077543427
$ws.Range("J7").Value = 0.9877039077543425
$ws.Range("M7").Value = 9.360041333333335
$ws.Range("O7").Value = 0.594282779407651
$ws.Range("P7").Value = 0.5942827794076508
$ws.Range("S7").Value = 0.5869754235320489
$ws.Range("T7").Value = 0.5869754235320486

# Row 8
$ws.Range("I8").Value = 0.01219848585420627
$ws.Range("J8").Value = 0.01219848585420626
$ws.Range("M8").Value = 3.442633333333333
$ws.Range("N8").Value = 10.3279
$ws.Range("O8").Value = 0.2185778494939793
$ws.Range("P8").Value = 0.2185778494939793
$ws.Range("Q8").Value = 5.634235516677778
$ws.Range("R8").Value = 50.7081196501
$ws.Range("S8").Value = 0.002666318805095133
$ws.Range("T8").Value = 0.002666318805095132

# Row 9
$ws.Range("I9").Value = 0.01219848585420627
$ws.Range("J9").Value = 0.01219848585420626
$ws.Range("O9").Value = 0.1871393710983698
$ws.Range("P9").Value = 0.1871393710983698
$ws.Range("S9").Value = 0.002282816971108521
$ws.Range("T9").Value = 0.00228281697110852

# Row 10
$ws.Range("I10").Value = 0.01219848585420627
$ws.Range("J10").Value = 0.01219848585420626
$ws.Range("M10").Value = 9.360041333333335
$ws.Range("O10").Value = 0.594282779407651
$ws.Range("P10").Value = 0.5942827794076508
$ws.Range("Q10").Value = 15.31870292639512
$ws.Range("S10").Value = 0.007249350078002614
$ws.Range("T10").Value = 0.00724935007800261

